$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 7.805543114114612
$ws.Cells.Item(2, 3).Value = 5.44426934738694
$ws.Cells.Item(2, 4).Value = 5.980106129458768
$ws.Cells.Item(2, 5).Value = 12.92265138345084
$ws.Cells.Item(2, 7).Value = 27.38410292719514
$ws.Cells.Item(2, 8).Value = 13.96336833258216
$ws.Cells.Item(2, 9).Value = 20.20127425907823
$ws.Cells.Item(2, 11).Value = 8.278630340669769
$ws.Cells.Item(2, 13).Value = 13.20517529870976
$ws.Cells.Item(2, 15).Value = 21.06785478488508
$ws.Cells.Item(3, 2).Value = 7.479159255367423
$ws.Cells.Item(3, 3).Value = 5.276570393855876
$ws.Cells.Item(3, 4).Value = 5.86015818202964
$ws.Cells.Item(3, 5).Value = 12.72239842815129
$ws.Cells.Item(3, 7).Value = 27.42748847833633
$ws.Cells.Item(3, 8).Value = 14.00926801709297
$ws.Cells.Item(3, 9).Value = 20.28952537903529
$ws.Cells.Item(3, 11).Value = 8.003866026527515
$ws.Cells.Item(3, 13).Value = 13.02095456359948
$ws.Cells.Item(3, 15).Value = 21.13700669877316
$ws.Cells.Item(4, 2).Value = 7.27202708196739
$ws.Cells.Item(4, 3).Value = 5.169923302844849
$ws.Cells.Item(4, 4).Value = 5.786987836945771
$ws.Cells.Item(4, 5).Value = 12.60251652323681
$ws.Cells.Item(4, 7).Value = 27.46358855484575
$ws.Cells.Item(4, 8).Value = 14.03975236640628
$ws.Cells.Item(4, 9).Value = 20.34774840779203
$ws.Cells.Item(4, 11).Value = 7.828683050233664
$ws.Cells.Item(4, 13).Value = 12.90957888070791
$ws.Cells.Item(4, 15).Value = 21.18421335048026
$ws.Cells.Item(5, 2).Value = 7.186046953095028
$ws.Cells.Item(5, 3).Value = 5.125579076165792
$ws.Cells.Item(5, 4).Value = 5.757337262259002
$ws.Cells.Item(5, 5).Value = 12.55450536322625
$ws.Cells.Item(5, 7).Value = 27.48066925665675
$ws.Cells.Item(5, 8).Value = 14.05275354206671
$ws.Cells.Item(5, 9).Value = 20.37248871154386
$ws.Cells.Item(5, 11).Value = 7.755730820374788
$ws.Cells.Item(5, 13).Value = 12.8646847351423
$ws.Cells.Item(5, 15).Value = 21.20464120862175
$ws.Cells.Item(6, 2).Value = 7.171678891226096
$ws.Cells.Item(6, 3).Value = 5.118163577986689
$ws.Cells.Item(6, 4).Value = 5.752425264320433
$ws.Cells.Item(6, 5).Value = 12.54658595310631
$ws.Cells.Item(6, 7).Value = 27.48364830972525
$ws.Cells.Item(6, 8).Value = 14.05494731372553
$ws.Cells.Item(6, 9).Value = 20.37665802954271
$ws.Cells.Item(6, 11).Value = 7.743524741875031
$ws.Cells.Item(6, 13).Value = 12.85726144450951
$ws.Cells.Item(6, 15).Value = 21.20810507030126
$ws.Cells.Item(7, 2).Value = 7.2708737187096
$ws.Cells.Item(7, 3).Value = 5.169328786766433
$ws.Cells.Item(7, 4).Value = 5.78658722090662
$ws.Cells.Item(7, 5).Value = 12.60186553060271
$ws.Cells.Item(7, 7).Value = 27.46380933021215
$ws.Cells.Item(7, 8).Value = 14.0399253625837
$ws.Cells.Item(7, 9).Value = 20.34807796014449
$ws.Cells.Item(7, 11).Value = 7.827705431778571
$ws.Cells.Item(7, 13).Value = 12.90897135608934
$ws.Cells.Item(7, 15).Value = 21.18448403033309
$ws.Cells.Item(8, 2).Value = 7.694470670097245
$ws.Cells.Item(8, 3).Value = 5.387233502249907
$ws.Cells.Item(8, 4).Value = 5.938676644876238
$ws.Cells.Item(8, 5).Value = 12.85300780484579
$ws.Cells.Item(8, 7).Value = 27.39709333522215
$ws.Cells.Item(8, 8).Value = 13.97871667175683
$ws.Cells.Item(8, 9).Value = 20.23086483306856
$ws.Cells.Item(8, 11).Value = 8.185272614216451
$ws.Cells.Item(8, 13).Value = 13.14132916507122
$ws.Cells.Item(8, 15).Value = 21.09071137484074
$ws.Cells.Item(9, 2).Value = 8.467181329125125
$ws.Cells.Item(9, 3).Value = 5.783762989600422
$ws.Cells.Item(9, 4).Value = 6.238706399734584
$ws.Cells.Item(9, 5).Value = 13.36680354881786
$ws.Cells.Item(9, 7).Value = 27.34169096407835
$ws.Cells.Item(9, 8).Value = 13.87696062738328
$ws.Cells.Item(9, 9).Value = 20.03308547907721
$ws.Cells.Item(9, 11).Value = 8.832564087607906
$ws.Cells.Item(9, 13).Value = 13.60829015189443
$ws.Cells.Item(9, 15).Value = 20.94461658633934
$ws.Cells.Item(10, 2).Value = 8.99440114842624
$ws.Cells.Item(10, 3).Value = 6.054457610591172
$ws.Cells.Item(10, 4).Value = 6.457618987744307
$ws.Cells.Item(10, 5).Value = 13.75310396539264
$ws.Cells.Item(10, 7).Value = 27.3473592170473
$ws.Cells.Item(10, 8).Value = 13.81335509482726
$ws.Cells.Item(10, 9).Value = 19.90739627071385
$ws.Cells.Item(10, 11).Value = 9.272358921661176
$ws.Cells.Item(10, 13).Value = 13.95494867215105
$ws.Cells.Item(10, 15).Value = 20.86048791541589
$ws.Cells.Item(11, 2).Value = 9.224604169110467
$ws.Cells.Item(11, 3).Value = 6.172794813153521
$ws.Cells.Item(11, 4).Value = 6.556359879665584
$ws.Cells.Item(11, 5).Value = 13.92987728697771
$ws.Cells.Item(11, 7).Value = 27.36005671854768
$ws.Cells.Item(11, 8).Value = 13.78684451749192
$ws.Cells.Item(11, 9).Value = 19.85449178839747
$ws.Cells.Item(11, 11).Value = 9.464165391677207
$ws.Cells.Item(11, 13).Value = 14.1127384999112
$ws.Cells.Item(11, 15).Value = 20.82728810742974
$ws.Cells.Item(12, 2).Value = 9.310333641775316
$ws.Cells.Item(12, 3).Value = 6.216891998197392
$ws.Cells.Item(12, 4).Value = 6.593588835587793
$ws.Cells.Item(12, 5).Value = 13.9968952283855
$ws.Cells.Item(12, 7).Value = 27.36632127607113
$ws.Cells.Item(12, 8).Value = 13.77715449936689
$ws.Cells.Item(12, 9).Value = 19.83507392488242
$ws.Cells.Item(12, 11).Value = 9.535574900524313
$ws.Cells.Item(12, 13).Value = 14.17244650253387
$ws.Cells.Item(12, 15).Value = 20.81544788454895
$ws.Cells.Item(13, 2).Value = 9.291935298154462
$ws.Cells.Item(13, 3).Value = 6.207427010173711
$ws.Cells.Item(13, 4).Value = 6.585578746626576
$ws.Cells.Item(13, 5).Value = 13.98245938813167
$ws.Cells.Item(13, 7).Value = 27.36490732255732
$ws.Cells.Item(13, 8).Value = 13.77922589150195
$ws.Cells.Item(13, 9).Value = 19.83922849070701
$ws.Cells.Item(13, 11).Value = 9.520250520540738
$ws.Cells.Item(13, 13).Value = 14.15959014403203
$ws.Cells.Item(13, 15).Value = 20.8179652997623
$ws.Cells.Item(14, 2).Value = 9.231686415571271
$ws.Cells.Item(14, 3).Value = 6.176437162361319
$ws.Cells.Item(14, 4).Value = 6.559426201800875
$ws.Cells.Item(14, 5).Value = 13.93538974253739
$ws.Cells.Item(14, 7).Value = 27.36054292073647
$ws.Cells.Item(14, 8).Value = 13.7860403178235
$ws.Cells.Item(14, 9).Value = 19.85288191604782
$ws.Cells.Item(14, 11).Value = 9.470065003693419
$ws.Cells.Item(14, 13).Value = 14.11765190725066
$ws.Cells.Item(14, 15).Value = 20.82629932368916
$ws.Cells.Item(15, 2).Value = 9.194592679152361
$ws.Cells.Item(15, 3).Value = 6.157361270740102
$ws.Cells.Item(15, 4).Value = 6.543384735334066
$ws.Cells.Item(15, 5).Value = 13.9065662359826
$ws.Cells.Item(15, 7).Value = 27.35805925557181
$ws.Cells.Item(15, 8).Value = 13.79025980793068
$ws.Cells.Item(15, 9).Value = 19.86132529115704
$ws.Cells.Item(15, 11).Value = 9.439164611405856
$ws.Cells.Item(15, 13).Value = 14.09195618720085
$ws.Cells.Item(15, 15).Value = 20.83149953512732
$ws.Cells.Item(16, 2).Value = 8.979157827188009
$ws.Cells.Item(16, 3).Value = 6.046625240237065
$ws.Cells.Item(16, 4).Value = 6.451145525544167
$ws.Cells.Item(16, 5).Value = 13.74156618117666
$ws.Cells.Item(16, 7).Value = 27.34673316755686
$ws.Cells.Item(16, 8).Value = 13.81513643345352
$ws.Cells.Item(16, 9).Value = 19.91093980472261
$ws.Cells.Item(16, 11).Value = 9.259654562879989
$ws.Cells.Item(16, 13).Value = 13.94463370409903
$ws.Cells.Item(16, 15).Value = 20.86275986811426
$ws.Cells.Item(17, 2).Value = 8.844484148505058
$ws.Cells.Item(17, 3).Value = 5.977444188112586
$ws.Cells.Item(17, 4).Value = 6.394313965248119
$ws.Cells.Item(17, 5).Value = 13.64055915447432
$ws.Cells.Item(17, 7).Value = 27.34237810268015
$ws.Cells.Item(17, 8).Value = 13.83101854797991
$ws.Cells.Item(17, 9).Value = 19.94247202453225
$ws.Cells.Item(17, 11).Value = 9.147388753405432
$ws.Cells.Item(17, 13).Value = 13.85423896592729
$ws.Cells.Item(17, 15).Value = 20.88323775087654
$ws.Cells.Item(18, 2).Value = 8.766119686108716
$ws.Cells.Item(18, 3).Value = 5.937202220498486
$ws.Cells.Item(18, 4).Value = 6.361548853831338
$ws.Cells.Item(18, 5).Value = 13.58256500460204
$ws.Cells.Item(18, 7).Value = 27.34082567144519
$ws.Cells.Item(18, 8).Value = 13.84038164475525
$ws.Cells.Item(18, 9).Value = 19.96101054368188
$ws.Cells.Item(18, 11).Value = 9.08204094187532
$ws.Cells.Item(18, 13).Value = 13.8022583892396
$ws.Cells.Item(18, 15).Value = 20.89549322107805
$ws.Cells.Item(19, 2).Value = 8.739433536311871
$ws.Cells.Item(19, 3).Value = 5.923500320447845
$ws.Cells.Item(19, 4).Value = 6.35044318810236
$ws.Cells.Item(19, 5).Value = 13.56294900790404
$ws.Cells.Item(19, 7).Value = 27.3404635702322
$ws.Cells.Item(19, 8).Value = 13.84359099655765
$ws.Cells.Item(19, 9).Value = 19.96735636733831
$ws.Cells.Item(19, 11).Value = 9.059783275514585
$ws.Cells.Item(19, 13).Value = 13.78466249603168
$ws.Cells.Item(19, 15).Value = 20.8997245773213
$ws.Cells.Item(20, 2).Value = 8.858914374611491
$ws.Cells.Item(20, 3).Value = 5.98485547294237
$ws.Cells.Item(20, 4).Value = 6.400372064274936
$ws.Cells.Item(20, 5).Value = 13.6513014549222
$ws.Cells.Item(20, 7).Value = 27.34274311367189
$ws.Cells.Item(20, 8).Value = 13.82930425675534
$ws.Cells.Item(20, 9).Value = 19.93907374691167
$ws.Cells.Item(20, 11).Value = 9.159420175209879
$ws.Cells.Item(20, 13).Value = 13.86386080695445
$ws.Cells.Item(20, 15).Value = 20.8810084439786
$ws.Cells.Item(21, 2).Value = 9.249422568123897
$ws.Cells.Item(21, 3).Value = 6.185559200675304
$ws.Cells.Item(21, 4).Value = 6.567112550276621
$ws.Cells.Item(21, 5).Value = 13.94921369512872
$ws.Cells.Item(21, 7).Value = 27.36178532736527
$ws.Cells.Item(21, 8).Value = 13.78402928166171
$ws.Cells.Item(21, 9).Value = 19.84885484553464
$ws.Cells.Item(21, 11).Value = 9.484839173481681
$ws.Cells.Item(21, 13).Value = 14.12997181281574
$ws.Cells.Item(21, 15).Value = 20.8238315369864
$ws.Cells.Item(22, 2).Value = 9.496207848735471
$ws.Cells.Item(22, 3).Value = 6.312556636895521
$ws.Cells.Item(22, 4).Value = 6.675125144630837
$ws.Cells.Item(22, 5).Value = 14.14433759338919
$ws.Cells.Item(22, 7).Value = 27.38271818007934
$ws.Cells.Item(22, 8).Value = 13.75647367297267
$ws.Cells.Item(22, 9).Value = 19.79348235986028
$ws.Cells.Item(22, 11).Value = 9.690374677800712
$ws.Cells.Item(22, 13).Value = 14.30361148438072
$ws.Cells.Item(22, 15).Value = 20.79072986842926
$ws.Cells.Item(23, 2).Value = 9.365282698055335
$ws.Cells.Item(23, 3).Value = 6.245164809244238
$ws.Cells.Item(23, 4).Value = 6.617577671919009
$ws.Cells.Item(23, 5).Value = 14.04018110673869
$ws.Cells.Item(23, 7).Value = 27.37076936611
$ws.Cells.Item(23, 8).Value = 13.77099435626933
$ws.Cells.Item(23, 9).Value = 19.82270661089182
$ws.Cells.Item(23, 11).Value = 9.58134093074551
$ws.Cells.Item(23, 13).Value = 14.2109806684726
$ws.Cells.Item(23, 15).Value = 20.80800560919046
$ws.Cells.Item(24, 2).Value = 8.852393388935203
$ws.Cells.Item(24, 3).Value = 5.981506289652393
$ws.Cells.Item(24, 4).Value = 6.397633480808531
$ws.Cells.Item(24, 5).Value = 13.64644461821474
$ws.Cells.Item(24, 7).Value = 27.34257512888834
$ws.Cells.Item(24, 8).Value = 13.83007856436551
$ws.Cells.Item(24, 9).Value = 19.94060883049265
$ws.Cells.Item(24, 11).Value = 9.153983273142943
$ws.Cells.Item(24, 13).Value = 13.85951080388611
$ws.Cells.Item(24, 15).Value = 20.88201481096666
$ws.Cells.Item(25, 2).Value = 8.264928402334611
$ws.Cells.Item(25, 3).Value = 5.67999311305666
$ws.Cells.Item(25, 4).Value = 6.157624738612919
$ws.Cells.Item(25, 5).Value = 13.22595850899659
$ws.Cells.Item(25, 7).Value = 27.34855715768315
$ws.Cells.Item(25, 8).Value = 13.90253049274675
$ws.Cells.Item(25, 9).Value = 20.08314952260203
$ws.Cells.Item(25, 11).Value = 8.663548428964054
$ws.Cells.Item(25, 13).Value = 13.48109760273602
$ws.Cells.Item(25, 15).Value = 20.98007527258492
